{"js": "const replacements = [\n  [\"2024-04-25 Thursday\", \"2024-04-26 Friday\"],\n  [\"216\u00f78=27, 0\", \"728\u00f72=364, 0\"],\n  [\"740\u00f78=92, 4\", \"268\u00f77=38, 2\"],\n  [\"681\u00f78=85, 1\", \"568\u00f79=63, 1\"],\n  [\"484\u00f79=53, 7\", \"503\u00f76=83, 5\"],\n  [\"858\u00f73=286, 0\", \"645\u00f74=161, 1\"],\n  [\"342\u00f78=42, 6\", \"237\u00f76=39, 3\"],\n  [\"620\u00f73=206, 2\", \"580\u00f78=72, 4\"],\n  [\"568\u00f74=142, 0\", \"870\u00f72=435, 0\"],\n  [\"780\u00f78=97, 4\", \"820\u00f72=410, 0\"],\n  [\"581\u00f78=72, 5\", \"943\u00f75=188, 3\"],\n  [\"469\u00f72=234, 1\", \"955\u00f74=238, 3\"],\n  [\"897\u00f78=112, 1\", \"591\u00f79=65, 6\"],\n  [\"104\u00f76=17, 2\", \"457\u00f76=76, 1\"],\n  [\"674\u00f76=112, 2\", \"608\u00f74=152, 0\"],\n  [\"758\u00f73=252, 2\", \"148\u00f76=24, 4\"],\n  [\"259\u00f72=129, 1\", \"408\u00f77=58, 2\"],\n  [\"809\u00f76=134, 5\", \"756\u00f77=108, 0\"],\n  [\"373\u00f72=186, 1\", \"314\u00f75=62, 4\"],\n  [\"486\u00f73=162, 0\", \"456\u00f79=50, 6\"],\n  [\"823\u00f72=411, 1\", \"841\u00f74=210, 1\"],\n  [\"695\u00f78=86, 7\", \"326\u00f73=108, 2\"],\n  [\"575\u00f77=82, 1\", \"138\u00f74=34, 2\"],\n  [\"677\u00f78=84, 5\", \"995\u00f75=199, 0\"],\n  [\"366\u00f77=52, 2\", \"489\u00f76=81, 3\"],\n  [\"898\u00f76=149, 4\", \"891\u00f78=111, 3\"],\n];\n\nconst body = context.document.body;\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-04-25 Thursday', '2024-04-26 Friday'),\n    @('216\u00f78=27, 0', '728\u00f72=364, 0'),\n    @('740\u00f78=92, 4', '268\u00f77=38, 2'),\n    @('681\u00f78=85, 1', '568\u00f79=63, 1'),\n    @('484\u00f79=53, 7', '503\u00f76=83, 5'),\n    @('858\u00f73=286, 0', '645\u00f74=161, 1'),\n    @('342\u00f78=42, 6', '237\u00f76=39, 3'),\n    @('620\u00f73=206, 2', '580\u00f78=72, 4'),\n    @('568\u00f74=142, 0', '870\u00f72=435, 0'),\n    @('780\u00f78=97, 4', '820\u00f72=410, 0'),\n    @('581\u00f78=72, 5', '943\u00f75=188, 3'),\n    @('469\u00f72=234, 1', '955\u00f74=238, 3'),\n    @('897\u00f78=112, 1', '591\u00f79=65, 6'),\n    @('104\u00f76=17, 2', '457\u00f76=76, 1'),\n    @('674\u00f76=112, 2', '608\u00f74=152, 0'),\n    @('758\u00f73=252, 2', '148\u00f76=24, 4'),\n    @('259\u00f72=129, 1', '408\u00f77=58, 2'),\n    @('809\u00f76=134, 5', '756\u00f77=108, 0'),\n    @('373\u00f72=186, 1', '314\u00f75=62, 4'),\n    @('486\u00f73=162, 0', '456\u00f79=50, 6'),\n    @('823\u00f72=411, 1', '841\u00f74=210, 1'),\n    @('695\u00f78=86, 7', '326\u00f73=108, 2'),\n    @('575\u00f77=82, 1', '138\u00f74=34, 2'),\n    @('677\u00f78=84, 5', '995\u00f75=199, 0'),\n    @('366\u00f77=52, 2', '489\u00f76=81, 3'),\n    @('898\u00f76=149, 4', '891\u00f78=111, 3'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute([ref]$pair[0], [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$pair[1], [ref]2) | Out-Null\n}\n"}
